$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Append 36 new translation rows (377-412) for the "setup" / "driptip" / "mod" lab features,
# mirroring the existing import-sheet layout: col A = language code ("cs"),
# col B = translation key, col C = translation value.
# Copy row 376 formatting down first (keeps the "import" cell style), then fill in values.
$ws.Range("A376:C376").Copy() | Out-Null
$ws.Range("A377:C412").PasteSpecial(-4122) | Out-Null

$ws.Range("A377").Value = "cs"
$ws.Range("B377").Value = "lab.setup.menu"
$ws.Range("C377").Value = "Setupy"
$ws.Range("A378").Value = "cs"
$ws.Range("B378").Value = "lab.setup.title"
$ws.Range("C378").Value = "Setupy"
$ws.Range("A379").Value = "cs"
$ws.Range("B379").Value = "lab.setup.subtitle"
$ws.Range("C379").Value = "Setup je složení fyzických zařízení použitých pro vapování."
$ws.Range("A380").Value = "cs"
$ws.Range("B380").Value = "lab.setup.button.create"
$ws.Range("C380").Value = "Nový setup"
$ws.Range("A381").Value = "cs"
$ws.Range("B381").Value = "lab.setup.button.list"
$ws.Range("C381").Value = "Seznam setupů"
$ws.Range("A382").Value = "cs"
$ws.Range("B382").Value = "lab.setup.create.title"
$ws.Range("C382").Value = "Nový setup"
$ws.Range("A383").Value = "cs"
$ws.Range("B383").Value = "lab.setup.create.subtitle"
$ws.Range("C383").Value = "Setup je poslední součást, která je potřebná pro sledování požitků z vapingu."
$ws.Range("A384").Value = "cs"
$ws.Range("B384").Value = "lab.setup.create.submit"
$ws.Range("C384").Value = "Vytvořit setup"
$ws.Range("A385").Value = "cs"
$ws.Range("B385").Value = "lab.setup.name.label"
$ws.Range("C385").Value = "Název setupu"
$ws.Range("A386").Value = "cs"
$ws.Range("B386").Value = "lab.setup.description.label"
$ws.Range("C386").Value = "Popis"
$ws.Range("A387").Value = "cs"
$ws.Range("B387").Value = "lab.setup.driptipId.label"
$ws.Range("C387").Value = "Náústek"
$ws.Range("A388").Value = "cs"
$ws.Range("B388").Value = "lab.setup.buildId.label"
$ws.Range("C388").Value = "Build"
$ws.Range("A389").Value = "cs"
$ws.Range("B389").Value = "lab.setup.modId.label"
$ws.Range("C389").Value = "Mod"
$ws.Range("A390").Value = "cs"
$ws.Range("B390").Value = "lab.driptip.tooltip.create"
$ws.Range("C390").Value = "Vytvořit náústek"
$ws.Range("A391").Value = "cs"
$ws.Range("B391").Value = "lab.driptip.create.title"
$ws.Range("C391").Value = "Nový náústek"
$ws.Range("A392").Value = "cs"
$ws.Range("B392").Value = "lab.driptip.create.subtitle"
$ws.Range("C392").Value = "Nezdá se to, ale náústky jsou také důležité."
$ws.Range("A393").Value = "cs"
$ws.Range("B393").Value = "lab.driptip.code.label"
$ws.Range("C393").Value = "Kód"
$ws.Range("A394").Value = "cs"
$ws.Range("B394").Value = "lab.driptip.vendorId.label"
$ws.Range("C394").Value = "Výrobce"
$ws.Range("A395").Value = "cs"
$ws.Range("B395").Value = "lab.driptip.create.submit"
$ws.Range("C395").Value = "Vytvořit náústek"
$ws.Range("A396").Value = "cs"
$ws.Range("B396").Value = "lab.driptip.created.message"
$ws.Range("C396").Value = "Náústek [{{data.code}}] byl uložen."
$ws.Range("A397").Value = "cs"
$ws.Range("B397").Value = "lab.build.tooltip.create"
$ws.Range("C397").Value = "Vytvořit build"
$ws.Range("A398").Value = "cs"
$ws.Range("B398").Value = "lab.mod.tooltip.create"
$ws.Range("C398").Value = "Vytvořit mod"
$ws.Range("A399").Value = "cs"
$ws.Range("B399").Value = "lab.mod.create.title"
$ws.Range("C399").Value = "Nový mod"
$ws.Range("A400").Value = "cs"
$ws.Range("B400").Value = "lab.mod.create.subtitle"
$ws.Range("C400").Value = "Mod obecně zastupuje zařízení, ze kterého lze vapovat."
$ws.Range("A401").Value = "cs"
$ws.Range("B401").Value = "lab.mod.name.label"
$ws.Range("C401").Value = "Název modu"
$ws.Range("A402").Value = "cs"
$ws.Range("B402").Value = "lab.mod.power.label"
$ws.Range("C402").Value = "Výkon (watty)"
$ws.Range("A403").Value = "cs"
$ws.Range("B403").Value = "lab.mod.vendorId.label"
$ws.Range("C403").Value = "Výrobce"
$ws.Range("A404").Value = "cs"
$ws.Range("B404").Value = "error.Duplicate entry [z_setup_name_unique] of [z_setup]."
$ws.Range("C404").Value = "Jméno tohoto setupu je již obsazené, použijte prosím jiné."
$ws.Range("A405").Value = "cs"
$ws.Range("B405").Value = "lab.setup.created.message"
$ws.Range("C405").Value = "Setup [{{data.name}}] byl uložen."
$ws.Range("A406").Value = "cs"
$ws.Range("B406").Value = "lab.setup.list.title"
$ws.Range("C406").Value = "Seznam setupů"
$ws.Range("A407").Value = "cs"
$ws.Range("B407").Value = "lab.setup.table.name"
$ws.Range("C407").Value = "Název"
$ws.Range("A408").Value = "cs"
$ws.Range("B408").Value = "lab.setup.table.driptip"
$ws.Range("C408").Value = "Náústek"
$ws.Range("A409").Value = "cs"
$ws.Range("B409").Value = "lab.setup.table.build"
$ws.Range("C409").Value = "Build"
$ws.Range("A410").Value = "cs"
$ws.Range("B410").Value = "lab.setup.table.mod"
$ws.Range("C410").Value = "Mod"
$ws.Range("A411").Value = "cs"
$ws.Range("B411").Value = "lab.build.inline.atomizer.tooltip"
$ws.Range("C411").Value = "Atomizér"
$ws.Range("A412").Value = "cs"
$ws.Range("B412").Value = "lab.build.inline.wraps.tooltip"
$ws.Range("C412").Value = "Počet otoček na spirálce"

$ws.Range("B401").Select() | Out-Null

